# Generate Report for Handoff
# - Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the six files that just finished handoff generation
#   (rows 7-12 on each sheet).
# - Sets the "Priority" column for those same rows to "ht" (handoff type)
#   on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
for ($r = 7; $r -le 12; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-19 02:21:21"
}

# zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority"
for ($r = 7; $r -le 12; $r++) {
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-19 02:21:15"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
}

# de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority"
for ($r = 7; $r -le 12; $r++) {
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-19 02:21:21"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
